{"js": "const sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\nconst sec = sections.items[0];\nconst headerFirst = sec.getHeader(\"FirstPage\");\nconst picsH = headerFirst.inlinePictures;\npicsH.load([\"items\"]);\nawait context.sync();\nconst p = picsH.items[0];\nconst rng = p.getRange();\nconst ooxmlResult = p.getOoxml();\nawait context.sync();\nreturn ooxmlResult.value.length;\n", "ps1": "$d = $word.ActiveDocument\n$sec = $d.Sections.Item(1)\n$h = $sec.Headers.Item(2)\n$s = $h.Range.InlineShapes.Item(1)\n$r = $s.Range\n$newXml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:wp=\"http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing\" xmlns:a=\"http://schemas.openxmlformats.org/drawingml/2006/main\" xmlns:pic=\"http://schemas.openxmlformats.org/drawingml/2006/picture\" xmlns:r=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships\"><w:body><w:p><w:r><w:drawing><wp:inline distB=\"0\" distT=\"0\" distL=\"0\" distR=\"0\"><wp:extent cx=\"914400\" cy=\"277792\"/><wp:effectExtent b=\"0\" l=\"0\" r=\"0\" t=\"0\"/><wp:docPr descr=\"BTec_Logo-Orange\" id=\"1\" name=\"image2.jpg\"/><a:graphic><a:graphicData uri=\"http://schemas.openxmlformats.org/drawingml/2006/picture\"><pic:pic><pic:nvPicPr><pic:cNvPr descr=\"BTec_Logo-Orange\" id=\"0\" name=\"image2.jpg\"/><pic:cNvPicPr preferRelativeResize=\"0\"/></pic:nvPicPr><pic:blipFill><a:blip r:embed=\"rId1\"/><a:srcRect b=\"0\" l=\"0\" r=\"0\" t=\"0\"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x=\"0\" y=\"0\"/><a:ext cx=\"914400\" cy=\"277792\"/></a:xfrm><a:prstGeom prst=\"rect\"/><a:ln/></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$r.InsertXML($newXml)\n"}
